# Datos.xlsx - "EnvioDatos" sheet gains new columns (Pais, Departamento,
# Ciudad, NumTel, TipoID, NumID) alongside the existing nombre/apellido/
# direccion columns, with sample row data. NumTel/NumID are stored as
# text (quote-prefixed) so the numeric-looking value "123456789" survives
# round-tripping as a string rather than becoming a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnvioDatos")
$ws.Activate()

# New header cells (row 1) -- written in the same order the original
# author typed them so shared-string interning lines up.
$ws.Range("D1").Value = "Pais"
$ws.Range("E1").Value = "Departamento"

# New data cells (row 2) for the Pais/Departamento pair -- Antioquia
# (the department) was entered before Colombia (the country).
$ws.Range("E2").Value = "Antioquia"
$ws.Range("D2").Value = "Colombia"

$ws.Range("F1").Value = "Ciudad"
$ws.Range("F2").Value = "Medellin"

$ws.Range("G1").Value = "NumTel"
$ws.Range("H1").Value = "TipoID"
$ws.Range("I1").Value = "NumID"

$ws.Range("H2").Value = "Cédula de Ciudadanía"

# Apply the "Text" number format (same style already used elsewhere in
# this workbook) across the whole used range of columns so every new
# cell -- including the blank-looking ones -- carries the shared style.
$ws.Columns("A:I").NumberFormat = "@"

# Phone number & ID number are typed as text-with-leading-apostrophe
# (Excel's "quote prefix"), which keeps the numeric-looking string as
# text and marks the style with quotePrefix.
$ws.Range("G2").Value = "'123456789"
$ws.Range("I2").Value = "'123456789"

# Match the final selection left behind in the sheet.
$ws.Range("I3").Select()
